# Apply the portfolio-weights update:
#  - Add header row ("Ticker" / "Weight")
#  - Populate the ticker list (column A) and weight values (column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tickers = @(
    "AMZN",
    "ADBE",
    "INTC",
    "BG",
    "GLTR",
    "MET",
    "NEE",
    "IBM",
    "OHI",
    "NKE",
    "AXP",
    "CACI",
    "CNC",
    "EMR",
    "LDOS",
    "MDT",
    "MPLX",
    "SO"
)

$weights = @(
    0.285744577706335,
    0.194255422293665,
    0.0500000000000001,
    0.0100000000000006,
    0.0100000000000001,
    0.0100000000000002,
    0.01,
    0.01,
    0.00999999999999994,
    0.05,
    0.01,
    0.01,
    0.00999999999999989,
    0.00999999999999952,
    0.0100000000000001,
    0.0100000000000002,
    0.15,
    0.15
)

# Header row
$ws.Range("A1").Value = "Ticker"
$ws.Range("B1").Value = "Weight"

# Data rows start at row 2
for ($i = 0; $i -lt $tickers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $tickers[$i]
    $ws.Cells.Item($row, 2).Value = $weights[$i]
}
